$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(1, 2).Value = 28.94093322753906
$ws.Cells.Item(1, 3).Value = 31.72183799743652
$ws.Cells.Item(1, 4).Value = 27.09791946411133
$ws.Cells.Item(1, 5).Value = 28.86097145080566
$ws.Cells.Item(2, 2).Value = 38.1544075012207
$ws.Cells.Item(2, 3).Value = 41.35848236083984
$ws.Cells.Item(2, 4).Value = 37.01863098144531
$ws.Cells.Item(2, 5).Value = 38.48972320556641
$ws.Cells.Item(3, 2).Value = 39.22373199462891
$ws.Cells.Item(3, 3).Value = 41.19161605834961
$ws.Cells.Item(3, 4).Value = 38.12992858886719
$ws.Cells.Item(3, 5).Value = 39.33806228637695
$ws.Cells.Item(4, 2).Value = 38.71438598632812
$ws.Cells.Item(4, 3).Value = 41.70479202270508
$ws.Cells.Item(4, 4).Value = 38.06488418579102
$ws.Cells.Item(4, 5).Value = 39.23088836669922
$ws.Cells.Item(5, 2).Value = 41.78150939941406
$ws.Cells.Item(5, 3).Value = 44.35115051269531
$ws.Cells.Item(5, 4).Value = 38.58174133300781
$ws.Cells.Item(5, 5).Value = 40.93862533569336
$ws.Cells.Item(6, 2).Value = 39.36916351318359
$ws.Cells.Item(6, 3).Value = 42.84889602661133
$ws.Cells.Item(6, 4).Value = 36.46199035644531
$ws.Cells.Item(6, 5).Value = 38.82323837280273
$ws.Cells.Item(7, 2).Value = 37.40185546875
$ws.Cells.Item(7, 3).Value = 42.38032913208008
$ws.Cells.Item(7, 4).Value = 39.09905242919922
$ws.Cells.Item(7, 5).Value = 39.17512512207031
$ws.Cells.Item(8, 2).Value = 33.43342208862305
$ws.Cells.Item(8, 3).Value = 34.10016632080078
$ws.Cells.Item(8, 4).Value = 32.573974609375
$ws.Cells.Item(8, 5).Value = 33.32404708862305
$ws.Cells.Item(9, 2).Value = 32.00667190551758
$ws.Cells.Item(9, 3).Value = 35.03984832763672
$ws.Cells.Item(9, 4).Value = 31.15790939331055
$ws.Cells.Item(9, 5).Value = 32.44304656982422
$ws.Cells.Item(10, 2).Value = 34.69312286376953
$ws.Cells.Item(10, 3).Value = 37.1483268737793
$ws.Cells.Item(10, 4).Value = 35.09734725952148
$ws.Cells.Item(10, 5).Value = 35.52105331420898
$ws.Cells.Item(11, 2).Value = 34.56444549560547
$ws.Cells.Item(11, 3).Value = 38.73139190673828
$ws.Cells.Item(11, 4).Value = 32.91348266601562
$ws.Cells.Item(11, 5).Value = 34.79392623901367
$ws.Cells.Item(12, 2).Value = 34.70332717895508
$ws.Cells.Item(12, 3).Value = 37.11955642700195
$ws.Cells.Item(12, 4).Value = 32.11446380615234
$ws.Cells.Item(12, 5).Value = 34.17467880249023
$ws.Cells.Item(13, 2).Value = 36.15987777709961
$ws.Cells.Item(13, 3).Value = 39.22965621948242
$ws.Cells.Item(13, 4).Value = 34.58393096923828
$ws.Cells.Item(13, 5).Value = 36.26144027709961
$ws.Cells.Item(14, 2).Value = 35.70352935791016
$ws.Cells.Item(14, 3).Value = 37.88385772705078
$ws.Cells.Item(14, 4).Value = 31.85197448730469
$ws.Cells.Item(14, 5).Value = 34.418701171875
$ws.Cells.Item(15, 2).Value = 38.61904144287109
$ws.Cells.Item(15, 3).Value = 40.88594436645508
$ws.Cells.Item(15, 4).Value = 35.36655807495117
$ws.Cells.Item(15, 5).Value = 37.69879913330078
$ws.Cells.Item(16, 2).Value = 40.76079940795898
$ws.Cells.Item(16, 3).Value = 43.88357162475586
$ws.Cells.Item(16, 4).Value = 39.00116348266602
$ws.Cells.Item(16, 5).Value = 40.77998352050781
$ws.Cells.Item(17, 2).Value = 39.04787445068359
$ws.Cells.Item(17, 3).Value = 42.90390396118164
$ws.Cells.Item(17, 4).Value = 38.77669525146484
$ws.Cells.Item(17, 5).Value = 39.88165664672852
$ws.Cells.Item(18, 2).Value = 36.29543304443359
$ws.Cells.Item(18, 3).Value = 40.55448913574219
$ws.Cells.Item(18, 4).Value = 36.54620742797852
$ws.Cells.Item(18, 5).Value = 37.41373443603516
$ws.Cells.Item(19, 2).Value = 36.64297400580512
$ws.Cells.Item(19, 3).Value = 39.61321205563016
$ws.Cells.Item(19, 4).Value = 35.24654748704698
$ws.Cells.Item(19, 5).Value = 36.75376118554009
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(1, 2).Value = 29.96093368530273
$ws.Cells.Item(1, 3).Value = 32.40523529052734
$ws.Cells.Item(1, 4).Value = 27.6295108795166
$ws.Cells.Item(1, 5).Value = 29.57317352294922
$ws.Cells.Item(2, 2).Value = 39.12648010253906
$ws.Cells.Item(2, 3).Value = 42.06500625610352
$ws.Cells.Item(2, 4).Value = 37.65140151977539
$ws.Cells.Item(2, 5).Value = 39.25460433959961
$ws.Cells.Item(3, 2).Value = 40.30384826660156
$ws.Cells.Item(3, 3).Value = 41.93997192382812
$ws.Cells.Item(3, 4).Value = 38.93620300292969
$ws.Cells.Item(3, 5).Value = 40.22319030761719
$ws.Cells.Item(4, 2).Value = 39.92569732666016
$ws.Cells.Item(4, 3).Value = 41.97026062011719
$ws.Cells.Item(4, 4).Value = 38.30780029296875
$ws.Cells.Item(4, 5).Value = 39.81715393066406
$ws.Cells.Item(5, 2).Value = 43.04160308837891
$ws.Cells.Item(5, 3).Value = 44.98368835449219
$ws.Cells.Item(5, 4).Value = 38.7957649230957
$ws.Cells.Item(5, 5).Value = 41.48055267333984
$ws.Cells.Item(6, 2).Value = 39.93154907226562
$ws.Cells.Item(6, 3).Value = 43.31013488769531
$ws.Cells.Item(6, 4).Value = 36.78837966918945
$ws.Cells.Item(6, 5).Value = 39.23550796508789
$ws.Cells.Item(7, 2).Value = 37.87023162841797
$ws.Cells.Item(7, 3).Value = 42.66083526611328
$ws.Cells.Item(7, 4).Value = 39.47819900512695
$ws.Cells.Item(7, 5).Value = 39.58274841308594
$ws.Cells.Item(8, 2).Value = 34.73881530761719
$ws.Cells.Item(8, 3).Value = 34.83522415161133
$ws.Cells.Item(8, 4).Value = 34.00452423095703
$ws.Cells.Item(8, 5).Value = 34.51005554199219
$ws.Cells.Item(9, 2).Value = 33.23057556152344
$ws.Cells.Item(9, 3).Value = 36.39508819580078
$ws.Cells.Item(9, 4).Value = 32.36769485473633
$ws.Cells.Item(9, 5).Value = 33.68436431884766
$ws.Cells.Item(10, 2).Value = 35.08209609985352
$ws.Cells.Item(10, 3).Value = 37.21687698364258
$ws.Cells.Item(10, 4).Value = 36.34161758422852
$ws.Cells.Item(10, 5).Value = 36.12432098388672
$ws.Cells.Item(11, 2).Value = 35.18029403686523
$ws.Cells.Item(11, 3).Value = 39.3519172668457
$ws.Cells.Item(11, 4).Value = 33.29235458374023
$ws.Cells.Item(11, 5).Value = 35.2869758605957
$ws.Cells.Item(12, 2).Value = 35.36347198486328
$ws.Cells.Item(12, 3).Value = 37.71631622314453
$ws.Cells.Item(12, 4).Value = 32.61042785644531
$ws.Cells.Item(12, 5).Value = 34.73553085327148
$ws.Cells.Item(13, 2).Value = 38.11630249023438
$ws.Cells.Item(13, 3).Value = 40.9373779296875
$ws.Cells.Item(13, 4).Value = 35.44962310791016
$ws.Cells.Item(13, 5).Value = 37.61105728149414
$ws.Cells.Item(14, 2).Value = 36.99066543579102
$ws.Cells.Item(14, 3).Value = 38.53049468994141
$ws.Cells.Item(14, 4).Value = 32.24367523193359
$ws.Cells.Item(14, 5).Value = 35.05496978759766
$ws.Cells.Item(15, 2).Value = 39.52053451538086
$ws.Cells.Item(15, 3).Value = 42.23884582519531
$ws.Cells.Item(15, 4).Value = 36.53905487060547
$ws.Cells.Item(15, 5).Value = 38.82441711425781
$ws.Cells.Item(16, 2).Value = 41.37936019897461
$ws.Cells.Item(16, 3).Value = 44.41862487792969
$ws.Cells.Item(16, 4).Value = 39.55954360961914
$ws.Cells.Item(16, 5).Value = 41.35459136962891
$ws.Cells.Item(17, 2).Value = 39.290771484375
$ws.Cells.Item(17, 3).Value = 43.25550079345703
$ws.Cells.Item(17, 4).Value = 39.14389038085938
$ws.Cells.Item(17, 5).Value = 40.19565582275391
$ws.Cells.Item(18, 2).Value = 36.98944854736328
$ws.Cells.Item(18, 3).Value = 41.38068771362305
$ws.Cells.Item(18, 4).Value = 37.08298110961914
$ws.Cells.Item(18, 5).Value = 38.0638542175293
$ws.Cells.Item(19, 2).Value = 37.55792660183377
$ws.Cells.Item(19, 3).Value = 40.31178262498644
$ws.Cells.Item(19, 4).Value = 35.90125815073649
$ws.Cells.Item(19, 5).Value = 37.47848468356662
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(1, 2).Value = 30.35678672790527
$ws.Cells.Item(1, 3).Value = 32.67437362670898
$ws.Cells.Item(1, 4).Value = 28.15261459350586
$ws.Cells.Item(1, 5).Value = 30.0123119354248
$ws.Cells.Item(2, 2).Value = 39.61936187744141
$ws.Cells.Item(2, 3).Value = 42.48798370361328
$ws.Cells.Item(2, 4).Value = 38.03919982910156
$ws.Cells.Item(2, 5).Value = 39.68436813354492
$ws.Cells.Item(3, 2).Value = 40.73922729492188
$ws.Cells.Item(3, 3).Value = 42.28229522705078
$ws.Cells.Item(3, 4).Value = 39.34452819824219
$ws.Cells.Item(3, 5).Value = 40.62541198730469
$ws.Cells.Item(4, 2).Value = 40.34889221191406
$ws.Cells.Item(4, 3).Value = 42.51802825927734
$ws.Cells.Item(4, 4).Value = 38.64346313476562
$ws.Cells.Item(4, 5).Value = 40.22359466552734
$ws.Cells.Item(5, 2).Value = 43.42312622070312
$ws.Cells.Item(5, 3).Value = 45.37884521484375
$ws.Cells.Item(5, 4).Value = 39.0314826965332
$ws.Cells.Item(5, 5).Value = 41.77330017089844
$ws.Cells.Item(6, 2).Value = 40.30735397338867
$ws.Cells.Item(6, 3).Value = 43.636962890625
$ws.Cells.Item(6, 4).Value = 37.04833602905273
$ws.Cells.Item(6, 5).Value = 39.53677749633789
$ws.Cells.Item(7, 2).Value = 38.21192932128906
$ws.Cells.Item(7, 3).Value = 42.98377227783203
$ws.Cells.Item(7, 4).Value = 39.72381210327148
$ws.Cells.Item(7, 5).Value = 39.88833618164062
$ws.Cells.Item(8, 2).Value = 35.38027572631836
$ws.Cells.Item(8, 3).Value = 35.30629730224609
$ws.Cells.Item(8, 4).Value = 34.50374984741211
$ws.Cells.Item(8, 5).Value = 35.04494094848633
$ws.Cells.Item(9, 2).Value = 33.89922714233398
$ws.Cells.Item(9, 3).Value = 37.08980941772461
$ws.Cells.Item(9, 4).Value = 32.94660186767578
$ws.Cells.Item(9, 5).Value = 34.31687545776367
$ws.Cells.Item(10, 2).Value = 35.56650924682617
$ws.Cells.Item(10, 3).Value = 37.503173828125
$ws.Cells.Item(10, 4).Value = 37.33901596069336
$ws.Cells.Item(10, 5).Value = 36.71079254150391
$ws.Cells.Item(11, 2).Value = 35.50518417358398
$ws.Cells.Item(11, 3).Value = 39.59118270874023
$ws.Cells.Item(11, 4).Value = 33.50859069824219
$ws.Cells.Item(11, 5).Value = 35.54305648803711
$ws.Cells.Item(12, 2).Value = 35.80410766601562
$ws.Cells.Item(12, 3).Value = 38.03202056884766
$ws.Cells.Item(12, 4).Value = 32.95048141479492
$ws.Cells.Item(12, 5).Value = 35.10023880004883
$ws.Cells.Item(13, 2).Value = 38.86020660400391
$ws.Cells.Item(13, 3).Value = 41.84979248046875
$ws.Cells.Item(13, 4).Value = 35.99428176879883
$ws.Cells.Item(13, 5).Value = 38.26963043212891
$ws.Cells.Item(14, 2).Value = 37.35730361938477
$ws.Cells.Item(14, 3).Value = 38.91839981079102
$ws.Cells.Item(14, 4).Value = 32.68972778320312
$ws.Cells.Item(14, 5).Value = 35.47522735595703
$ws.Cells.Item(15, 2).Value = 40.13093948364258
$ws.Cells.Item(15, 3).Value = 42.92842483520508
$ws.Cells.Item(15, 4).Value = 37.02784729003906
$ws.Cells.Item(15, 5).Value = 39.37753295898438
$ws.Cells.Item(16, 2).Value = 41.68775939941406
$ws.Cells.Item(16, 3).Value = 44.73407363891602
$ws.Cells.Item(16, 4).Value = 40.01743698120117
$ws.Cells.Item(16, 5).Value = 41.73886108398438
$ws.Cells.Item(17, 2).Value = 39.39449310302734
$ws.Cells.Item(17, 3).Value = 43.57820892333984
$ws.Cells.Item(17, 4).Value = 39.50307846069336
$ws.Cells.Item(17, 5).Value = 40.44224548339844
$ws.Cells.Item(18, 2).Value = 37.60507965087891
$ws.Cells.Item(18, 3).Value = 41.86933517456055
$ws.Cells.Item(18, 4).Value = 37.43199157714844
$ws.Cells.Item(18, 5).Value = 38.5466423034668
$ws.Cells.Item(19, 2).Value = 38.01098685794406
$ws.Cells.Item(19, 3).Value = 40.74238777160645
$ws.Cells.Item(19, 4).Value = 36.32756890190972
$ws.Cells.Item(19, 5).Value = 37.90611913469102
